$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Tidy up two grammar-checker split runs (merge the runs back into
#    a single run and drop the now-redundant proofErr markers). The
#    visible text is unchanged; only the run/proofErr structure is
#    simplified, exactly like Word does when it re-normalizes runs
#    after a Find & Replace over the spanned text.
# ------------------------------------------------------------------

$apostrophe = [char]0x2019

$r1 = $d.Content
$old1 = " stays unchanged at all times. Each of the three operations has a certain "
$null = $r1.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $old1, 2)

$r2 = $d.Content
$old2 = " operation is 0 if it doesn" + $apostrophe + "t actually change the character."
$null = $r2.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $old2, 2)

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark (the marker Word drops at the last
#    edit position) away from the "Minimum Edit Distance" heading and
#    onto the point where new text/solution content was added, in the
#    middle of "symbols" in the SymbolMultiplication problem
#    statement. Bookmark names are unique, so adding a new "_GoBack"
#    bookmark automatically removes the old one.
# ------------------------------------------------------------------

$r3 = $d.Content
$found = $r3.Find.Execute("all sym", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r3.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $r3)
}
